$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-05-18 Sunday" "2025-05-19 Monday"

Replace-Text "36×19=684" "18×43=774"
Replace-Text "64×46=2944" "73×18=1314"
Replace-Text "99×30=2970" "41×66=2706"
Replace-Text "60×62=3720" "94×31=2914"
Replace-Text "15×67=1005" "20×37=740"

Replace-Text "19×98=1862" "75×94=7050"
Replace-Text "92×35=3220" "67×74=4958"
Replace-Text "95×48=4560" "17×69=1173"
Replace-Text "65×30=1950" "67×40=2680"
Replace-Text "39×21=819" "74×97=7178"

Replace-Text "95×52=4940" "84×49=4116"
Replace-Text "23×88=2024" "23×78=1794"
Replace-Text "13×60=780" "65×62=4030"
Replace-Text "85×47=3995" "51×79=4029"
Replace-Text "61×13=793" "34×65=2210"

Replace-Text "71×16=1136" "72×28=2016"
Replace-Text "60×68=4080" "87×30=2610"
Replace-Text "18×75=1350" "62×44=2728"
Replace-Text "20×39=780" "20×95=1900"
Replace-Text "86×86=7396" "78×89=6942"

Replace-Text "21×91=1911" "81×42=3402"
Replace-Text "59×86=5074" "80×41=3280"
Replace-Text "96×26=2496" "82×25=2050"
Replace-Text "25×44=1100" "47×98=4606"
Replace-Text "43×23=989" "87×79=6873"
